$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list: refreshed Price (D) and Volume(1h) (E) values per row.
# For Price cells, force text format so purely-numeric-looking strings (e.g. "582.61")
# are stored as text (matching the source inline-string data) instead of being
# auto-converted by Excel into floating point numbers.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.989.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.34%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.465.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.88%  "
# Row 4
$ws.Range("E4").Value = "  +0.05%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.56%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.71%  "
# Row 7
$ws.Range("E7").Value = "  +0.11%  "
# Row 8
$ws.Range("E8").Value = "  -2.27%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.465.42"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.87%  "
# Row 10
$ws.Range("E10").Value = "  -2.71%  "
# Row 11
$ws.Range("E11").Value = "  -0.17%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.12%  "
# Row 13
$ws.Range("E13").Value = "  -4.05%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.67%  "
# Row 15
$ws.Range("E15").Value = "  -1.90%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.863.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.38%  "
# Row 17
$ws.Range("E17").Value = "  -4.79%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.474.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.43%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.60%  "
# Row 20
$ws.Range("E20").Value = "  -5.10%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.70%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.46%  "
# Row 23
$ws.Range("E23").Value = "  -0.05%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.70%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.44%  "
# Row 26
$ws.Range("E26").Value = "  -7.27%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.52%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -55.54%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.566.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.80%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "513.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.04%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0901"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.16%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -9.04%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.28%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.68%  "
# Row 35
$ws.Range("E35").Value = "  +0.12%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.117"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.15%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.04%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.35"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.35%  "
# Row 40
$ws.Range("E40").Value = "  -6.62%  "
# Row 41
$ws.Range("E41").Value = "  -0.09%  "
# Row 42
$ws.Range("E42").Value = "  -6.65%  "
# Row 43
$ws.Range("E43").Value = "  -7.03%  "
# Row 44
$ws.Range("E44").Value = "  -7.38%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.60%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.73"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.15%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "140.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.04%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.14%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.514"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.94%  "
# Row 50
$ws.Range("E50").Value = "  -12.50%  "
# Row 51
$ws.Range("E51").Value = "  -7.48%  "
